# Update "想去人数" (want-to-go count) figures in column F across the four
# sheets of the workbook, matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 454
$ws1.Range("F5").Value = 8746
$ws1.Range("F7").Value = 11138
$ws1.Range("F8").Value = 93
$ws1.Range("F13").Value = 123
$ws1.Range("F14").Value = 226
$ws1.Range("F15").Value = 299
$ws1.Range("F17").Value = 32
$ws1.Range("F20").Value = 421
$ws1.Range("F22").Value = 1895
$ws1.Range("F23").Value = 707
$ws1.Range("F24").Value = 633
$ws1.Range("F25").Value = 358
$ws1.Range("F28").Value = 605
$ws1.Range("F29").Value = 57
$ws1.Range("F30").Value = 1293
$ws1.Range("F31").Value = 26
$ws1.Range("F32").Value = 10
$ws1.Range("F35").Value = 1426
$ws1.Range("F37").Value = 5
$ws1.Range("F39").Value = 302
$ws1.Range("F40").Value = 33
$ws1.Range("F41").Value = 143
$ws1.Range("F43").Value = 382
$ws1.Range("F47").Value = 31
$ws1.Range("F48").Value = 154
$ws1.Range("F49").Value = 142

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 55
$ws2.Range("F14").Value = 30
$ws2.Range("F15").Value = 12
$ws2.Range("F16").Value = 13

# 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2837
$ws3.Range("F4").Value = 349

# 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 454
$ws4.Range("F8").Value = 8746
$ws4.Range("F10").Value = 11138
$ws4.Range("F11").Value = 93
$ws4.Range("F14").Value = 123
$ws4.Range("F15").Value = 299
$ws4.Range("F16").Value = 32
$ws4.Range("F19").Value = 1895
$ws4.Range("F20").Value = 707
$ws4.Range("F21").Value = 633
$ws4.Range("F22").Value = 358
$ws4.Range("F25").Value = 605
$ws4.Range("F26").Value = 55
$ws4.Range("F29").Value = 1293
$ws4.Range("F30").Value = 26
$ws4.Range("F31").Value = 10
$ws4.Range("F33").Value = 30
$ws4.Range("F34").Value = 12
$ws4.Range("F35").Value = 13
$ws4.Range("F37").Value = 1426
$ws4.Range("F42").Value = 382
$ws4.Range("F47").Value = 31
$ws4.Range("F48").Value = 154
$ws4.Range("F49").Value = 142
